$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: VENTAS POR GRUPO
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Update M137 (GUERRERO FAREZ FABIAN MAURICIO / WONG SANCHEZ CLAUDIA PAULINA)
$ws1.Range("M137").Value = 794.88

# Insert a new row at 247 for a new client "DANIELA ELIZABETH BECERRA BECERRA"
# under advisor "OFICINA-CATAECSA", pushing all following rows down by one.
$ws1.Rows.Item(247).Insert()
$ws1.Range("A247").Value = "OFICINA-CATAECSA"
$ws1.Range("B247").Value = "DANIELA ELIZABETH BECERRA BECERRA"
$ws1.Range("C247:R247").Value = 0

# Update the trailing "N de 288" -> "N de 289" counter row, now shifted to row 291
$ws1.Range("C291").Value = "5 de 289"
$ws1.Range("D291").Value = "16 de 289"
$ws1.Range("E291").Value = "10 de 289"
$ws1.Range("F291").Value = "1 de 289"
$ws1.Range("G291").Value = "0 de 289"
$ws1.Range("H291").Value = "9 de 289"
$ws1.Range("I291").Value = "16 de 289"
$ws1.Range("J291").Value = "1 de 289"
$ws1.Range("K291").Value = "1 de 289"
$ws1.Range("L291").Value = "24 de 289"
$ws1.Range("M291").Value = "62 de 289"
$ws1.Range("N291").Value = "3 de 289"
$ws1.Range("O291").Value = "3 de 289"
$ws1.Range("P291").Value = "4 de 289"
$ws1.Range("Q291").Value = "5 de 289"
$ws1.Range("R291").Value = "0 de 289"

# ---------------------------------------------------------------------------
# Sheet 2: VENTA MENSUAL
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Update F137 (GUERRERO FAREZ FABIAN MAURICIO / WONG SANCHEZ CLAUDIA PAULINA)
$ws2.Range("F137").Value = 794.88

# Insert the same new row at 247
$ws2.Rows.Item(247).Insert()
$ws2.Range("A247").Value = "OFICINA-CATAECSA"
$ws2.Range("B247").Value = "DANIELA ELIZABETH BECERRA BECERRA"
$ws2.Range("C247:G247").Value = 0

# Update the grand-total row, now shifted to row 291 (only GRANITO/F changes)
$ws2.Range("F291").Value = 258218.96

# ---------------------------------------------------------------------------
# Sheet 3: CUMPLIMIENTO MENSUAL
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# GUERRERO FAREZ FABIAN MAURICIO / PORCELANATO
$ws3.Range("D50").Value = 58960.61
$ws3.Range("E50").Value = -7134.150000000001
$ws3.Range("F50").Value = 1.137654588023184

# TOTAL row
$ws3.Range("D123").Value = 284420.9
$ws3.Range("E123").Value = 188348.4947370426
$ws3.Range("F123").Value = 0.6016059905024027
